# Appends the new check_availability log entries (rows 10-24).
# Mirrors the diff: dimension A1:F9 -> A1:F24, 15 new rows appended.
# Row 16 records a Selenium stack-trace failure in column D instead of the
# usual "No availability..." message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = '2024-09-09 16:22:02'
$ws.Range("B10").Value = 'check_availability'
$ws.Range("C10").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D10").Value = 'No availability for the selected date.'
$ws.Range("E10").Value = "'2024-09-09"
$ws.Range("F10").Value = '16:22:02'

$ws.Range("A11").Value = '2024-09-09 16:22:34'
$ws.Range("B11").Value = 'check_availability'
$ws.Range("C11").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D11").Value = 'No availability for the selected date.'
$ws.Range("E11").Value = "'2024-09-09"
$ws.Range("F11").Value = '16:22:34'

$ws.Range("A12").Value = '2024-09-09 16:23:05'
$ws.Range("B12").Value = 'check_availability'
$ws.Range("C12").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D12").Value = 'No availability for the selected date.'
$ws.Range("E12").Value = "'2024-09-09"
$ws.Range("F12").Value = '16:23:05'

$ws.Range("A13").Value = '2024-09-09 16:23:36'
$ws.Range("B13").Value = 'check_availability'
$ws.Range("C13").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D13").Value = 'No availability for the selected date.'
$ws.Range("E13").Value = "'2024-09-09"
$ws.Range("F13").Value = '16:23:36'

$ws.Range("A14").Value = '2024-09-09 16:24:07'
$ws.Range("B14").Value = 'check_availability'
$ws.Range("C14").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D14").Value = 'No availability for the selected date.'
$ws.Range("E14").Value = "'2024-09-09"
$ws.Range("F14").Value = '16:24:07'

$ws.Range("A15").Value = '2024-09-09 16:24:39'
$ws.Range("B15").Value = 'check_availability'
$ws.Range("C15").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D15").Value = 'No availability for the selected date.'
$ws.Range("E15").Value = "'2024-09-09"
$ws.Range("F15").Value = '16:24:39'

$ws.Range("A16").Value = '2024-09-09 16:26:19'
$ws.Range("B16").Value = 'check_availability'
$ws.Range("C16").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$d16 = @'
Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label-wrapper button[aria-label*='30']"}
  (Session info: chrome=128.0.6613.120); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF7F082B5D2+29090]
	(No symbol) [0x00007FF7F079E689]
	(No symbol) [0x00007FF7F065B1CA]
	(No symbol) [0x00007FF7F06AEFD7]
	(No symbol) [0x00007FF7F06AF22C]
	(No symbol) [0x00007FF7F06F97F7]
	(No symbol) [0x00007FF7F06D672F]
	(No symbol) [0x00007FF7F06F65D9]
	(No symbol) [0x00007FF7F06D6493]
	(No symbol) [0x00007FF7F06A09B1]
	(No symbol) [0x00007FF7F06A1B11]
	GetHandleVerifier [0x00007FF7F0B48C5D+3295277]
	GetHandleVerifier [0x00007FF7F0B94843+3605523]
	GetHandleVerifier [0x00007FF7F0B8A707+3564247]
	GetHandleVerifier [0x00007FF7F08E6EB6+797318]
	(No symbol) [0x00007FF7F07A980F]
	(No symbol) [0x00007FF7F07A53F4]
	(No symbol) [0x00007FF7F07A5580]
	(No symbol) [0x00007FF7F0794A1F]
	BaseThreadInitThunk [0x00007FFC979C257D+29]
	RtlUserThreadStart [0x00007FFC9896AF28+40]

'@
$ws.Range("D16").Value = $d16
$ws.Range("E16").Value = "'2024-09-09"
$ws.Range("F16").Value = '16:26:19'
# The multi-line stack trace makes Excel auto-expand the row height; reset it
# to the sheet default so row 16 matches the others (no custom height stored).
$ws.Rows.Item(16).AutoFit()

$ws.Range("A17").Value = '2024-09-09 16:33:12'
$ws.Range("B17").Value = 'check_availability'
$ws.Range("C17").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D17").Value = 'No availability for the selected date.'
$ws.Range("E17").Value = "'2024-09-09"
$ws.Range("F17").Value = '16:33:12'

$ws.Range("A18").Value = '2024-09-09 16:33:44'
$ws.Range("B18").Value = 'check_availability'
$ws.Range("C18").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D18").Value = 'No availability for the selected date.'
$ws.Range("E18").Value = "'2024-09-09"
$ws.Range("F18").Value = '16:33:44'

$ws.Range("A19").Value = '2024-09-09 16:43:45'
$ws.Range("B19").Value = 'check_availability'
$ws.Range("C19").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D19").Value = 'No availability for the selected date.'
$ws.Range("E19").Value = "'2024-09-09"
$ws.Range("F19").Value = '16:43:45'

$ws.Range("A20").Value = '2024-09-09 16:44:17'
$ws.Range("B20").Value = 'check_availability'
$ws.Range("C20").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D20").Value = 'No availability for the selected date.'
$ws.Range("E20").Value = "'2024-09-09"
$ws.Range("F20").Value = '16:44:17'

$ws.Range("A21").Value = '2024-09-09 16:51:07'
$ws.Range("B21").Value = 'check_availability'
$ws.Range("C21").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D21").Value = 'No availability for the selected date.'
$ws.Range("E21").Value = "'2024-09-09"
$ws.Range("F21").Value = '16:51:07'

$ws.Range("A22").Value = '2024-09-09 16:51:33'
$ws.Range("B22").Value = 'check_availability'
$ws.Range("C22").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D22").Value = 'No availability for the selected date.'
$ws.Range("E22").Value = "'2024-09-09"
$ws.Range("F22").Value = '16:51:33'

$ws.Range("A23").Value = '2024-09-09 16:51:54'
$ws.Range("B23").Value = 'check_availability'
$ws.Range("C23").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D23").Value = 'No availability for the selected date.'
$ws.Range("E23").Value = "'2024-09-09"
$ws.Range("F23").Value = '16:51:54'

$ws.Range("A24").Value = '2024-09-09 16:52:35'
$ws.Range("B24").Value = 'check_availability'
$ws.Range("C24").Value = 'https://www.opentable.com/r/bar-spero-washington/'
$ws.Range("D24").Value = 'No availability for the selected date.'
$ws.Range("E24").Value = "'2024-09-09"
$ws.Range("F24").Value = '16:52:35'

